# Fruta / hortaliza, semanal
# Insert a new weekly record as row 48, pushing existing rows 48:83 down to 49:84.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 48 (shifts rows 48:83 -> 49:84)
$ws.Rows("48:48").Insert()

# Populate the new row 48 with the new weekly record
$ws.Range("A48").Value = 11
$ws.Range("B48").Value = "Vega Monumental Concepción"
$ws.Range("C48").Value = "Bíobío"
$ws.Range("D48").Value = 44489
$ws.Range("E48").Value = 8
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100102
$ws.Range("H48").Value = "Cítricos"
$ws.Range("I48").Value = 100102004
$ws.Range("J48").Value = "Mandarina"
$ws.Range("K48").Value = "Murcott"
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 650
$ws.Range("N48").Value = 7000
$ws.Range("O48").Value = 7500
$ws.Range("P48").Value = 7231
$ws.Range("Q48").Value = "$/caja 18 kilos"
$ws.Range("R48").Value = "Región de O'Higgins"
$ws.Range("S48").Value = 402
$ws.Range("T48").Value = 18
